$wb = $excel.ActiveWorkbook

# 1. Update the comment on the "Valuation" sheet, cell A55, explaining the
#    6-month finance input (construction period note).
$wsVal = $wb.Worksheets.Item("Valuation")
$comment = $wsVal.Range("A55").Comment
$comment.Text("Carter, Jo:" + [char]10 + "If construction period is 12 months, enter 6 in this cell")

# 2. Replace the named valuer list on the "Dropdown lists" sheet with
#    anonymised placeholder names (NAME1..NAME8), merging/removing the
#    real names that used to live in the shared strings table.
$wsDrop = $wb.Worksheets.Item("Dropdown lists")
$names = @("NAME1", "NAME2", "NAME3", "NAME4", "NAME5", "NAME6", "NAME7", "NAME8")
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 18 + $i
    $wsDrop.Range("A$row").Value = $names[$i]
}

# 3. Restore/update the view selections: "Dropdown lists" now has A26
#    selected, while "Valuation" scrolls back to the top with B3 selected
#    and remains the active (tabSelected) sheet.
$wsDrop.Activate()
$wsDrop.Range("A26").Select()

$wsVal.Activate()
$wsVal.Range("B3").Select()
